#
# Applies the "Created geo map with datapoints" commit's textual edits:
#  - "Model used:" paragraph -> "Techniques" + " used:" (two runs), followed
#    by new "Preprocessing" / "Cleaning" / "Modeling ..." bullet paragraphs
#    and two blank paragraphs. The _GoBack bookmark moves from the "Hard"
#    paragraph onto the new "Preprocessing" paragraph.
#  - Several proofErr (grammar-check) markers get added/shifted around
#    existing runs further down in the "Easy" / code-walkthrough section.
#
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: locate a paragraph (by exact visible text, paragraph mark
# trimmed) anywhere in the active document. Re-scans every call so it
# stays correct as earlier edits shift paragraph indices around.
# ---------------------------------------------------------------------
function Get-ParaByText {
    param([string]$targetText)
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $d.Paragraphs($i)
        $t = $para.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $targetText) {
            return $para
        }
    }
    return $null
}

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# 1) Drop the _GoBack bookmark currently sitting after "Hard" — it will
#    be re-created further up, after the new "Preprocessing" line.
# ---------------------------------------------------------------------
$hardPara = Get-ParaByText("Hard")
$hardXml = $pkgHeader + '<w:p><w:r><w:t>Hard</w:t></w:r></w:p>' + $pkgFooter
$hardPara.Range.InsertXML($hardXml)

# ---------------------------------------------------------------------
# 2) Replace "Model used:" with "Techniques"/" used:" plus the new
#    Preprocessing / Cleaning / Modeling bullets and two blank lines.
#    (InsertXML drops the very last wholly-empty <w:p> in a payload, so
#    one extra trailing <w:p/> is added to compensate and land on two
#    real blank paragraphs.)
# ---------------------------------------------------------------------
$modelPara = Get-ParaByText("Model used:")
$modelBody = '<w:p><w:r><w:t>Techniques</w:t></w:r><w:r><w:t xml:space="preserve"> used:</w:t></w:r></w:p>' + `
    '<w:p><w:r><w:tab/><w:t>Preprocessing</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' + `
    '<w:p><w:r><w:tab/><w:t>Cleaning</w:t></w:r></w:p>' + `
    '<w:p><w:r><w:tab/><w:t>Modeling \u2013 Regression, Decision Trees, Random Forest</w:t></w:r></w:p>' + `
    '<w:p/><w:p/><w:p/>'
$modelBody = $modelBody.Replace('\u2013', [char]0x2013)
$modelXml = $pkgHeader + $modelBody + $pkgFooter
$modelPara.Range.InsertXML($modelXml)

# ---------------------------------------------------------------------
# 3) "X  = time" -> gramStart/"X  ="/gramEnd/" time"
# ---------------------------------------------------------------------
$xPara = Get-ParaByText("X  = time")
$xBody = '<w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr>' + `
    '<w:proofErr w:type="gramStart"/><w:r><w:t>X  =</w:t></w:r><w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> time</w:t></w:r></w:p>'
$xXml = $pkgHeader + $xBody + $pkgFooter
$xPara.Range.InsertXML($xXml)

# ---------------------------------------------------------------------
# 4) "Xtest, ytest = trin test split(x,y, test,size0.2, random)"
#    -> " test " / gramStart / "split(" / spellStart+gramEnd / "x,y" / ...
# ---------------------------------------------------------------------
$splitPara = Get-ParaByText("Xtest, ytest = trin test split(x,y, test,size0.2, random)")
$splitBody = '<w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Xtest</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">, </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>ytest</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> = </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>trin</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> test </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/><w:r><w:t>split(</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t>x,y</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>, test,size0.2, random)</w:t></w:r></w:p>'
$splitXml = $pkgHeader + $splitBody + $pkgFooter
$splitPara.Range.InsertXML($splitXml)

# ---------------------------------------------------------------------
# 5) "Reg = linear regression().fit(xtrain,ytrain)"
#    -> "Reg = linear regression(" / gramStart / ").fit" / gramEnd / "(" / ...
# ---------------------------------------------------------------------
$regPara = Get-ParaByText("Reg = linear regression().fit(xtrain,ytrain)")
$regBody = '<w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr>' + `
    '<w:r><w:t>Reg = linear regression(</w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/><w:r><w:t>).fit</w:t></w:r><w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t>(</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>xtrain,ytrain</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>)</w:t></w:r></w:p>'
$regXml = $pkgHeader + $regBody + $pkgFooter
$regPara.Range.InsertXML($regXml)

# ---------------------------------------------------------------------
# 6) "Pred = reg.predict(X_test)" -> wrap reg.predict in gramStart/gramEnd
#    (in addition to its existing spellStart/spellEnd).
# ---------------------------------------------------------------------
$predPara = Get-ParaByText("Pred = reg.predict(X_test)")
$predBody = '<w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Pred</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> = </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>reg.predict</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t>(</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>X_test</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>)</w:t></w:r></w:p>'
$predXml = $pkgHeader + $predBody + $pkgFooter
$predPara.Range.InsertXML($predXml)
